$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo in existing row 17 (Utami et al) Details cell: "bene- fits" -> "benefits" ---
$ws.Range("E17").Value = "in this study we aimed to: (a) build an approach capable of assessing the cost, relative benefits and cost-effectiveness of implementing threat management strategies that improve a broad range of values in multifunctional areas; (b) bring together and build key information to help managers and stakeholders understand the values, goals, threats, total management costs and opportunities for achieving goals for values, using the TNBB as a case study; and (c) deliver a set of costed, prioritized strategies for achieving goals across multiple important values of the TNBB."

# --- Add missing Subject cell to existing row 18 (Cullen) ---
# Copy formatting (style) from D19, an existing Subject-column cell with the correct wrap/font style
$ws.Range("D19").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D18").Value = "Reviews SCP literature"

# --- Append new rows 20-26 with paper details ---
# Source cells used as format donors:
#   A17/C17 -> plain (non-wrapped) style used by Authors/year columns
#   B17/D17/E17 -> wrapped style used by Title/Subject/Details columns
#   E2        -> the distinct "coloured" wrapped style used once for D25

# Row 20
$ws.Range("A17").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A20").Value = "Vivitskaia et al"
$ws.Range("B17").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B20").Value = "Linking threat maps with management to guide conservation investment"
$ws.Range("A17").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C20").Value = 2020
$ws.Range("B17").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D20").Value = "Built impact maps for coastal waters around the world and identified best strategies for investement for top 10 places"
$ws.Range("B17").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E20").Value = "We rebuild cumulative impact maps by stressor type (climate change, marine and land) at a global scale to evaluate the expected effectiveness of various management strategies for all coastal territories. Key disparities were found between broad-scale management of marine ecosystems and the dominant stressors, with existing management in tropical island nations likely insufficient to address intense impacts from climate change. These countries also typically had low performance on governance indicators, suggesting challenges in implementing new mitigation. We highlight trade-offs in making decisions for stressor mitigation and offer strategic guidance on identifying locations to target management of marine, land, or climate impacts."
$ws.Rows.Item(20).RowHeight = 84.6

# Row 21
$ws.Range("A17").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A21").Value = "Jepson et al"
$ws.Range("B17").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B21").Value = "Protected area asset stewardship"
$ws.Range("A17").Copy()
$ws.Range("C21").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C21").Value = 2017
$ws.Range("B17").Copy()
$ws.Range("D21").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D21").Value = "Argue for framing PAs as spatial assets for generating investment"
$ws.Range("B17").Copy()
$ws.Range("E21").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E21").Value = "Our asset framework offers a complementary investment approach and proposition. The world is awash with capital but typical returns on investments are historically low. This, in combination with ideas of impact investing, is generating a demand for conservation investment products. Once formalised, a PA asset approach would create the capacity to optimize PA assets (as sites or networks) in terms of their spatial location, investment profile, and the forms of value they generate over time. It would also support and extend the programme of work on PA management effectiveness (Coad et al., 2015) through providing a framework to assess the social, economic and cultural benefits of PAs."
$ws.Rows.Item(21).RowHeight = 84.6

# Row 22
$ws.Range("A17").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A22").Value = "Kearney et al"
$ws.Range("B17").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B22").Value = "Estimating the benefit of well-managed protected areas for threatened speices conservation"
$ws.Range("A17").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C22").Value = 2018
$ws.Range("B17").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D22").Value = "Use Australia PA network as example of under-resourced PAs not protecting species"
$ws.Range("B17").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E22").Value = "Good reference for the case for adequately investing in protected areas. Also warn against expanding PA networks without ensureing adequte resources for those PAs"
$ws.Rows.Item(22).RowHeight = 48.6

# Row 23
$ws.Range("A17").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A23").Value = "Coad et al "
$ws.Range("B17").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B23").Value = "Widespread shortfalls in protected area resourcing undermine efforts to conserve biodiversity "
$ws.Range("A17").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C23").Value = 2019
$ws.Range("B17").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D23").Value = "assess >2100 PA reports and quantify how many are under funded and resourced"
$ws.Rows.Item(23).RowHeight = 48.6

# Row 24
$ws.Range("A17").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A24").Value = "Pringle"
$ws.Range("B17").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B24").Value = "Upgrading protected areas to conserve wild biodiversity"
$ws.Range("A17").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C24").Value = 2017
$ws.Range("B17").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D24").Value = "Uses two case studies to show how PAs can be expanded and rewilded"
$ws.Range("B17").Copy()
$ws.Range("E24").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E24").Value = "International agreements mandate the expansion of Earth's protected-area network as a bulwark against the continued extinction of wild populations, species, and ecosystems. Yet many protected areas are underfunded, poorly managed, and ecologically damaged; the conundrum is how to increase their coverage and effectiveness simultaneously. Worldwide, enormous potential for biodiversity conservation can be realized by upgrading existing nature reserves while harmonizing them with the needs and aspirations of their constituencies."
$ws.Rows.Item(24).RowHeight = 60.6

# Row 25
$ws.Range("A17").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A25").Value = "Robinson et al"
$ws.Range("B17").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B25").Value = "Incorporating land tenure security into conservation"
$ws.Range("A17").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C25").Value = 2017
$ws.Range("E2").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D25").Value = "present a framework that identifies three common ways in which land tenure security can impact human and conservation outcomes"
$ws.Range("B17").Copy()
$ws.Range("E25").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E25").Value = "Insecure land tenure plagues many developing and tropical regions, often where conservation concerns are highest. We present a framework that identifies three common ways in which land tenure security can impact human and conservation outcomes, and suggest practical ways to distill tenure and tenure security issues for a given location. We conclude with steps for considering tenure security issues in the context of conservation projects and identify areas for future research."
$ws.Rows.Item(25).RowHeight = 72.6

# Row 26
$ws.Range("A17").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A26").Value = "de Oliveira et al "
$ws.Range("B17").Copy()
$ws.Range("B26").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B26").Value = "The financial needs vs. the realities of in situ conservation: an analysis of federal funding for protected areas in Brazil's Caatinga"
$ws.Range("A17").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C26").Value = 2017
$ws.Rows.Item(26).RowHeight = 36.6

# --- Update the selected cell to match the final workbook state ---
$ws.Range("E26").Select()

Write-Output "done"
